$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" header in F1, copying the header formatting from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (2-24)
$timestamps = @(
    "2021-10-05 13:41:55.457618",
    "2021-10-05 13:41:55.457630",
    "2021-10-05 13:41:55.457633",
    "2021-10-05 13:41:55.457636",
    "2021-10-05 13:41:55.457640",
    "2021-10-05 13:41:55.457643",
    "2021-10-05 13:41:55.457646",
    "2021-10-05 13:41:55.457649",
    "2021-10-05 13:41:55.457652",
    "2021-10-05 13:41:55.457655",
    "2021-10-05 13:41:55.457658",
    "2021-10-05 13:41:55.457661",
    "2021-10-05 13:41:55.457664",
    "2021-10-05 13:41:55.457667",
    "2021-10-05 13:41:55.457670",
    "2021-10-05 13:41:55.457673",
    "2021-10-05 13:41:55.457676",
    "2021-10-05 13:41:55.457679",
    "2021-10-05 13:41:55.457682",
    "2021-10-05 13:41:55.457685",
    "2021-10-05 13:41:55.457688",
    "2021-10-05 13:41:55.457691",
    "2021-10-05 13:41:55.457694"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
